$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.393753
$ws.Range("N2").Value = 1.181259
$ws.Range("O2").Value = 0.05217225506577142
$ws.Range("P2").Value = 0.05217225506577142
$ws.Range("Q2").Value = 0.067186336892
$ws.Range("R2").Value = 0.6046770320279999
$ws.Range("S2").Value = 0.05217225506577142
$ws.Range("T2").Value = 0.05217225506577142

# Row 3 updates
$ws.Range("O3").Value = 0.4686908567124721
$ws.Range("P3").Value = 0.4686908567124721
$ws.Range("S3").Value = 0.4686908567124721
$ws.Range("T3").Value = 0.4686908567124721

# Row 4 updates
$ws.Range("M4").Value = 3.616128666666667
$ws.Range("N4").Value = 10.848386
$ws.Range("O4").Value = 0.4791368882217565
$ws.Range("P4").Value = 0.4791368882217565
$ws.Range("Q4").Value = 0.6170224451457778
$ws.Range("R4").Value = 5.553202006312
$ws.Range("S4").Value = 0.4791368882217565
$ws.Range("T4").Value = 0.4791368882217565
